$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new "LinElast" data block
$ws.Range("A9").Value = "X-Elem"
$ws.Range("B9").Value = "Y-Elem"
$ws.Range("C9").Value = "Temp"
$ws.Range("D9").Value = "Step Diff"

# Row 10 - first data point, no step diff (nothing to diff against)
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 38.5108

# Row 11 - plain (non-shared) difference formula
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 20
$ws.Range("C11").Value = 40.345
$ws.Range("D11").Formula = "=C11-C10"

# Row 12
$ws.Range("A12").Value = 16
$ws.Range("B12").Value = 40
$ws.Range("C12").Value = 42.1826
$ws.Range("D12").Formula = "=C12-C11"

# Row 13
$ws.Range("A13").Value = 32
$ws.Range("B13").Value = 80
$ws.Range("C13").Value = 44.021
$ws.Range("D13").Formula = "=C13-C12"

# Row 14
$ws.Range("A14").Value = 64
$ws.Range("B14").Value = 160
$ws.Range("C14").Value = 45.8595
$ws.Range("D14").Formula = "=C14-C13"

# Row 15
$ws.Range("A15").Value = 96
$ws.Range("B15").Value = 240
$ws.Range("C15").Value = 46.935
$ws.Range("D15").Formula = "=C15-C14"

# Update the active selection to match the saved workbook state
$ws.Range("F17").Select()
